$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New hyperlink row (A25): "Proposé par timothée rapin : ... donné"
# linking to the forum thread, with the hyperlink's "display" text using the
# "... donné par timothée rapin" word order (matches author's original paste).
$targetUrl = "https://codes-sources.commentcamarche.net/forum/affich-371867-plein-ecran-c-console"
$cellText  = "Proposé par timothée rapin : $targetUrl donné"
$displayText = "$targetUrl donné par timothée rapin"

# Create the hyperlink first (this both writes the cell text and applies the
# "Lien hypertexte" style), then overwrite the visible text and re-apply the
# exact same named cell style so it matches the other rows in the sheet.
$ws.Hyperlinks.Add($ws.Range("A25"), $targetUrl, [Type]::Missing, [Type]::Missing, $displayText)
$ws.Range("A25").Value = $cellText
$ws.Range("A25").Style = "Lien hypertexte"

# Mirror the author's new selection rectangle (A27:A28).
$ws.Range("A27:A28").Select()
